# Updates cryptos list figures (prices / 1h volume %) to match the
# latest scrape, and rotates several rows 36-51 to reflect coins that
# entered/left the tracked ranking (names, links and values shift by one
# position; Frax drops out, Elrond is newly added at the bottom).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "30.560.16" }
    @{ Cell = "E2"; Value = "  -0.21%  " }
    @{ Cell = "D3"; Value = "1.914.64" }
    @{ Cell = "E3"; Value = "  -0.50%  " }
    @{ Cell = "E4"; Value = "  +0.03%  " }
    @{ Cell = "D5"; Value = "'244.70" }
    @{ Cell = "E5"; Value = "  -0.94%  " }
    @{ Cell = "E6"; Value = "  -0.02%  " }
    @{ Cell = "D7"; Value = "'0.4849" }
    @{ Cell = "E7"; Value = "  +2.12%  " }
    @{ Cell = "D8"; Value = "'0.2888" }
    @{ Cell = "E8"; Value = "  -0.04%  " }
    @{ Cell = "D9"; Value = "'0.06790" }
    @{ Cell = "E9"; Value = "  -0.51%  " }
    @{ Cell = "D10"; Value = "'111.38" }
    @{ Cell = "E10"; Value = "  +5.73%  " }
    @{ Cell = "D11"; Value = "'19.38" }
    @{ Cell = "D12"; Value = "1.917.30" }
    @{ Cell = "E12"; Value = "  -0.36%  " }
    @{ Cell = "E13"; Value = "  -1.35%  " }
    @{ Cell = "D14"; Value = "'5.402" }
    @{ Cell = "E14"; Value = "  +1.27%  " }
    @{ Cell = "D15"; Value = "'0.6708" }
    @{ Cell = "E15"; Value = "  +0.34%  " }
    @{ Cell = "D16"; Value = "'294.21" }
    @{ Cell = "E16"; Value = "  +1.38%  " }
    @{ Cell = "D17"; Value = "30.549.86" }
    @{ Cell = "E17"; Value = "  -0.23%  " }
    @{ Cell = "D18"; Value = "'13.01" }
    @{ Cell = "E18"; Value = "  +0.42%  " }
    @{ Cell = "D19"; Value = "'0.000007597" }
    @{ Cell = "E19"; Value = "  -0.31%  " }
    @{ Cell = "D20"; Value = "'1.000" }
    @{ Cell = "E20"; Value = "  -0.02%  " }
    @{ Cell = "D21"; Value = "'5.526" }
    @{ Cell = "E21"; Value = "  -1.37%  " }
    @{ Cell = "D22"; Value = "2.165.46" }
    @{ Cell = "E22"; Value = "  -0.32%  " }
    @{ Cell = "E23"; Value = "  +0.01%  " }
    @{ Cell = "D24"; Value = "'6.420" }
    @{ Cell = "E24"; Value = "  -0.44%  " }
    @{ Cell = "D25"; Value = "'9.470" }
    @{ Cell = "E25"; Value = "  -0.14%  " }
    @{ Cell = "D26"; Value = "'166.15" }
    @{ Cell = "E26"; Value = "  -0.47%  " }
    @{ Cell = "D27"; Value = "'20.24" }
    @{ Cell = "E27"; Value = "  -4.59%  " }
    @{ Cell = "D28"; Value = "'2.078" }
    @{ Cell = "E28"; Value = "  -1.48%  " }
    @{ Cell = "D29"; Value = "'0.1064" }
    @{ Cell = "E29"; Value = "  -0.68%  " }
    @{ Cell = "D30"; Value = "'1.443" }
    @{ Cell = "E30"; Value = "  +2.94%  " }
    @{ Cell = "D31"; Value = "'4.123" }
    @{ Cell = "E31"; Value = "  -1.49%  " }
    @{ Cell = "D32"; Value = "'4.048" }
    @{ Cell = "E32"; Value = "  +0.17%  " }
    @{ Cell = "D33"; Value = "'0.04983" }
    @{ Cell = "E33"; Value = "  -0.86%  " }
    @{ Cell = "D34"; Value = "'0.7334" }
    @{ Cell = "E34"; Value = "  +0.26%  " }
    @{ Cell = "D35"; Value = "'1.141" }
    @{ Cell = "E35"; Value = "  -0.13%  " }
    @{ Cell = "B36"; Value = "HuobiToken" }
    @{ Cell = "C36"; Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht" }
    @{ Cell = "D36"; Value = "'2.711" }
    @{ Cell = "E36"; Value = "  -0.89%  " }
    @{ Cell = "B37"; Value = "VeChain" }
    @{ Cell = "C37"; Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet" }
    @{ Cell = "D37"; Value = "'0.02029" }
    @{ Cell = "E37"; Value = "  -1.57%  " }
    @{ Cell = "B38"; Value = "MXToken" }
    @{ Cell = "C38"; Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx" }
    @{ Cell = "D38"; Value = "'2.681" }
    @{ Cell = "E38"; Value = "  -0.29%  " }
    @{ Cell = "B39"; Value = "RenderToken" }
    @{ Cell = "C39"; Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr" }
    @{ Cell = "D39"; Value = "'2.020" }
    @{ Cell = "E39"; Value = "  -1.43%  " }
    @{ Cell = "B40"; Value = "Quant" }
    @{ Cell = "C40"; Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt" }
    @{ Cell = "D40"; Value = "'109.34" }
    @{ Cell = "E40"; Value = "  -2.10%  " }
    @{ Cell = "B41"; Value = "TheSandbox" }
    @{ Cell = "C41"; Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand" }
    @{ Cell = "D41"; Value = "'0.4435" }
    @{ Cell = "E41"; Value = "  +1.07%  " }
    @{ Cell = "B42"; Value = "TrustWalletToken" }
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt" }
    @{ Cell = "D42"; Value = "'0.8677" }
    @{ Cell = "E42"; Value = "  -0.80%  " }
    @{ Cell = "B43"; Value = "FraxShare" }
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs" }
    @{ Cell = "D43"; Value = "'5.833" }
    @{ Cell = "E43"; Value = "  -1.56%  " }
    @{ Cell = "B44"; Value = "PaxDollar" }
    @{ Cell = "C44"; Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp" }
    @{ Cell = "D44"; Value = "'1.000" }
    @{ Cell = "E44"; Value = "  -0.02%  " }
    @{ Cell = "B45"; Value = "Aave" }
    @{ Cell = "C45"; Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave" }
    @{ Cell = "D45"; Value = "'69.36" }
    @{ Cell = "E45"; Value = "  +2.45%  " }
    @{ Cell = "B46"; Value = "Aptos" }
    @{ Cell = "C46"; Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt" }
    @{ Cell = "D46"; Value = "'7.206" }
    @{ Cell = "E46"; Value = "  -1.12%  " }
    @{ Cell = "B47"; Value = "BitcoinSV" }
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv" }
    @{ Cell = "D47"; Value = "'48.52" }
    @{ Cell = "E47"; Value = "  +0.01%  " }
    @{ Cell = "B48"; Value = "EnergySwap" }
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" }
    @{ Cell = "D48"; Value = "'9.241" }
    @{ Cell = "E48"; Value = "  -0.72%  " }
    @{ Cell = "B49"; Value = "Algorand" }
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo" }
    @{ Cell = "D49"; Value = "'0.1226" }
    @{ Cell = "E49"; Value = "  -1.41%  " }
    @{ Cell = "B50"; Value = "WOONetwork" }
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo" }
    @{ Cell = "D50"; Value = "'0.2502" }
    @{ Cell = "E50"; Value = "  +0.52%  " }
    @{ Cell = "B51"; Value = "Elrond" }
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld" }
    @{ Cell = "D51"; Value = "'34.75" }
    @{ Cell = "E51"; Value = "  -0.54%  " }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

Write-Host ("Applied " + $updates.Count + " cell updates")
